# Plantillas modelos impuestos UHY.xlsx -- "working upload multiplefiles"
#
# On the "IVA" sheet, the July column (K) values were copied across into the
# August (L) and September (M) columns for a number of rows (the months had
# been left blank before).
#
# On the "IRPF 111" sheet, the June column (J) picked up values on two
# row groups (rows 8-10 and 20-22) that were previously blank.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "IVA": copy column K (julio) into L (agosto) and M (septiembre)
# for every row that already has a July figure.
# ---------------------------------------------------------------------
$wsIva = $wb.Worksheets.Item("IVA")

$ivaRows = @(8, 9, 10, 12, 20, 21, 22, 24, 31, 32, 35, 37, 41, 42, 45, 47, 57, 61, 62)

foreach ($r in $ivaRows) {
    $julyValue = $wsIva.Range("K$r").Value()
    $wsIva.Range("L$r").Value = $julyValue
    $wsIva.Range("M$r").Value = $julyValue
}

# ---------------------------------------------------------------------
# Sheet "IRPF 111": fill in the June (J) column figures.
# ---------------------------------------------------------------------
$wsIrpf111 = $wb.Worksheets.Item("IRPF 111")

$wsIrpf111.Range("J8").Value = 50
$wsIrpf111.Range("J9").Value = 91023.50999999999
$wsIrpf111.Range("J10").Value = 2825.92

$wsIrpf111.Range("J20").Value = 4
$wsIrpf111.Range("J21").Value = 1629.85
$wsIrpf111.Range("J22").Value = 128.3
